# Updated symbol list (cryptocurrency price refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.60"
$ws.Range("D3").Value = "'23.12"
$ws.Range("D4").Value = "'5.419"
$ws.Range("D5").Value = "'0.05940"
$ws.Range("D6").Value = "'3.454"
$ws.Range("D7").Value = "'6.530"
$ws.Range("D8").Value = "'0.8129"
$ws.Range("D9").Value = "'0.9094"
$ws.Range("D10").Value = "'0.1409"
$ws.Range("D11").Value = "'0.07484"
$ws.Range("D12").Value = "'0.03294"
$ws.Range("D13").Value = "'0.03060"
$ws.Range("D14").Value = "'0.09351"
$ws.Range("D15").Value = "'3.864"
$ws.Range("D16").Value = "'0.001574"
$ws.Range("D17").Value = "'0.04673"
$ws.Range("D18").Value = "'0.0005942"
$ws.Range("D19").Value = "'0.006066"
$ws.Range("D20").Value = "'0.004997"
$ws.Range("D21").Value = "'0.0009805"
$ws.Range("E22").Value = "21NitroExNTXBestin24h"
$ws.Range("D23").Value = "'3.612"
$ws.Range("D24").Value = "'2.146"
$ws.Range("D25").Value = "'0.3228"
$ws.Range("D27").Value = "'0.0002395"
$ws.Range("D40").Value = "'0.03941"
$ws.Range("D41").Value = "'0.006199"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").Value = "'0.1076"
$ws.Range("D43").Value = "'0.002621"
$ws.Range("D44").Value = "'0.007829"
$ws.Range("D45").Value = "'0.00005229"
$ws.Range("D48").Value = "'0.9004"
$ws.Range("D49").Value = "'0.002264"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("D51").Value = "'0.0002001"
